# Add newly-available rural production units to the Availability sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Availability")

$xlPasteFormats = -4122

# --- Row 16: new entry (SUPELCAMM03), styled like row 12's J cell (s=25) ---
$ws.Range("F15:I15").Copy() | Out-Null
$ws.Range("F16:I16").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J12").Copy() | Out-Null
$ws.Range("J16").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("F16:I16").Value = 1
$ws.Range("J16").Value = "SUPELCAMM03"

# --- Row 17: new entry (SUPH2ALKR1N), unstyled J cell (like row 15's J) ---
$ws.Range("F15:I15").Copy() | Out-Null
$ws.Range("F17:I17").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J15").Copy() | Out-Null
$ws.Range("J17").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("F17:I17").Value = 1
$ws.Range("J17").Value = "SUPH2ALKR1N"

# --- Row 18: new entry (SUPH2ALKR2N), styled like existing J18/L18 (s=27) ---
$ws.Range("F15:I15").Copy() | Out-Null
$ws.Range("F18:I18").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("F18:I18").Value = 1
$ws.Range("J18").Value = "SUPH2ALKR2N"

# --- Row 19: new entry (SUPH2PEMR1N), styled like existing J19/L19 (s=27) ---
$ws.Range("F15:I15").Copy() | Out-Null
$ws.Range("F19:I19").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("F19:I19").Value = 1
$ws.Range("J19").Value = "SUPH2PEMR1N"

# --- Row 20: new entry (SUPH2PEMR2N), unstyled J cell ---
$ws.Range("F15:I15").Copy() | Out-Null
$ws.Range("F20:I20").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J15").Copy() | Out-Null
$ws.Range("J20").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("F20:I20").Value = 1
$ws.Range("J20").Value = "SUPH2PEMR2N"

# --- Row 21: new entry (SUPH2SOER2N), unstyled J cell ---
$ws.Range("F15:I15").Copy() | Out-Null
$ws.Range("F21:I21").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("J15").Copy() | Out-Null
$ws.Range("J21").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("F21:I21").Value = 1
$ws.Range("J21").Value = "SUPH2SOER2N"

# --- Row 22: trailing blank row, J cell styled like row13/14's J cell (s=26) ---
$ws.Range("J13").Copy() | Out-Null
$ws.Range("J22").PasteSpecial($xlPasteFormats) | Out-Null

# --- Row-height touch-ups that accompany the new rural block ---
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 15
$ws.Rows.Item(13).AutoFit() | Out-Null
$ws.Rows.Item(14).AutoFit() | Out-Null
$ws.Rows.Item(15).RowHeight = 15

# Update selection to match the new active cell (cosmetic, matches the diff's sheetView)
$ws.Range("J22").Select() | Out-Null
